$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Componentes")

# Update the Column A unique ID values (text) per row as described in the diff.
$updates = @{
    2 = "6337025003091F0E"
    4 = "6337025103B2704E"
    6 = "6337025103B54D1F"
    8 = "6337025003015CBB"
    10 = "6337025003023A17"
    12 = "6337025103B1633F"
    14 = "633702500305697F"
    16 = "6337025103B95BE6"
    18 = "6337025003043055"
    20 = "63370250039940FA"
    22 = "6337025003A21C32"
    24 = "63370250030B6E40"
    26 = "63370250038F69E4"
    28 = "6337025003960DD5"
    30 = "63370250039D6C83"
    32 = "63370250030C743C"
    34 = "63370250030D241C"
    36 = "63370250030E3E36"
    38 = "6337025003135999"
    40 = "63370250031433D5"
    42 = "63370250031514E8"
    44 = "6337025003163F8C"
    46 = "63370250038772EA"
    48 = "6337025003074A57"
    50 = "63370250030840E8"
    52 = "63370250035F696A"
    54 = "6337025003655675"
    56 = "63370250038D7764"
    58 = "63370250039404AA"
    60 = "63370250032A76C5"
    61 = "63370250032B02A2"
    62 = "63370250032C001C"
    63 = "63370250032D576E"
    64 = "6337025003034189"
    66 = "6337025003062E85"
    68 = "6337025003173E73"
    71 = "63370250035E4ACE"
    73 = "6337025003614A1B"
    75 = "6337025003694BDA"
    77 = "63370250036E6B65"
    79 = "633702500373179D"
    81 = "6337025003784DE9"
    83 = "63370250037D1CDE"
    85 = "6337025003823D6D"
    87 = "63370250038C0A3D"
    89 = "63370250039366ED"
    91 = "6337025003AF093D"
    93 = "6337025003B03D9A"
    95 = "6337025103BB2B72"
    97 = "63370250039B58E9"
    99 = "633702500375649C"
    101 = "6337025003701FEA"
    103 = "63370250038472FA"
    105 = "63370250037F1D5D"
    107 = "63370250037A79B7"
    109 = "63370250036B4ADD"
    111 = "63370250036335FF"
    113 = "63370250030A1F3C"
    115 = "6337025003107934"
    117 = "633702500312078C"
    118 = "633702500311284F"
    120 = "63370250030F652B"
    122 = "6337025003A32AE1"
}

foreach ($row in $updates.Keys) {
    $ws.Cells.Item($row, 1).Value = $updates[$row]
}
